$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct "Marking" row (row 11): Right/Wrong marks
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Correct "Total" row (row 12): total marks scored and max marks text
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "80 / 112"
